$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlRight = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignRight

# Right-align the existing "Ngày" (date) column cells first, B4:B9.
# Doing this before touching row 10 reproduces the same cellXfs ordering
# Excel itself used: a plain "right aligned" style (no number format) gets
# created first, and a "right aligned + date number format" style second.
$ws.Range("B4:B9").HorizontalAlignment = $xlRight

# New week-2 row: task #2, its date, and the task description.
$ws.Range("A10").Value = 2
$ws.Range("B10").Value = Get-Date -Year 2023 -Month 1 -Day 10 -Hour 0 -Minute 0 -Second 0
$ws.Range("B10").HorizontalAlignment = $xlRight
$ws.Range("C10").Value = "Setup truyền thông LoRa"

# Move/restore the active selection the way the author's session ended up.
$ws.Range("B11").Select() | Out-Null
